$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1)
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F2").Value = 285
$wsExhibit.Range("F4").Value = 7798
$wsExhibit.Range("F5").Value = 5685
$wsExhibit.Range("F6").Value = 473
$wsExhibit.Range("F8").Value = 11
$wsExhibit.Range("F10").Value = 261
$wsExhibit.Range("F11").Value = 259
$wsExhibit.Range("F12").Value = 57

# Sheet "全部类型" (sheet4)
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 285
$wsAll.Range("F4").Value = 7798
$wsAll.Range("F5").Value = 5685
$wsAll.Range("F6").Value = 473
$wsAll.Range("F8").Value = 11
$wsAll.Range("F10").Value = 261
$wsAll.Range("F13").Value = 259
$wsAll.Range("F14").Value = 57
